$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.219803940680094
$ws.Cells.Item(2, 4).Value = 0.1058422808887993
$ws.Cells.Item(2, 5).Value = 0.13028179933346
$ws.Cells.Item(2, 6).Value = 2.878478165862106
$ws.Cells.Item(2, 7).Value = 2.368769996223534
$ws.Cells.Item(2, 8).Value = 1.696868061313182
$ws.Cells.Item(2, 9).Value = 3.961193417866752
$ws.Cells.Item(2, 10).Value = 0.2151472402880472
$ws.Cells.Item(3, 2).Value = 1.06999749926689
$ws.Cells.Item(3, 4).Value = 0.1057099897463161
$ws.Cells.Item(3, 5).Value = 0.129376248143025
$ws.Cells.Item(3, 6).Value = 2.750918434283136
$ws.Cells.Item(3, 7).Value = 2.222640816841505
$ws.Cells.Item(3, 8).Value = 1.636692919714051
$ws.Cells.Item(3, 9).Value = 3.481483530428136
$ws.Cells.Item(3, 10).Value = 0.208773148833032
$ws.Cells.Item(4, 2).Value = 0.9775436500523824
$ws.Cells.Item(4, 4).Value = 0.1056316089112226
$ws.Cells.Item(4, 5).Value = 0.1288919000434596
$ws.Cells.Item(4, 6).Value = 2.674709327376632
$ws.Cells.Item(4, 7).Value = 2.13468666062019
$ws.Cells.Item(4, 8).Value = 1.600961959854857
$ws.Cells.Item(4, 9).Value = 3.186307868167489
$ws.Cells.Item(4, 10).Value = 0.205012059613054
$ws.Cells.Item(5, 2).Value = 0.9397517899388959
$ws.Cells.Item(5, 4).Value = 0.1056003683078792
$ws.Cells.Item(5, 5).Value = 0.1287124574646867
$ws.Cells.Item(5, 6).Value = 2.64417696227224
$ws.Cells.Item(5, 7).Value = 2.099280920985308
$ws.Cells.Item(5, 8).Value = 1.586703403558658
$ws.Cells.Item(5, 9).Value = 3.065876933592335
$ws.Cells.Item(5, 10).Value = 0.2035173224339673
$ws.Cells.Item(6, 2).Value = 0.9334695415148531
$ws.Cells.Item(6, 4).Value = 0.1055952226743999
$ws.Cells.Item(6, 5).Value = 0.128683741779934
$ws.Cells.Item(6, 6).Value = 2.639138474748648
$ws.Cells.Item(6, 7).Value = 2.093427924020915
$ws.Cells.Item(6, 8).Value = 1.584353915396832
$ws.Cells.Item(6, 9).Value = 3.045871178739986
$ws.Cells.Item(6, 10).Value = 0.2032714010193857
$ws.Cells.Item(7, 2).Value = 0.9770344427025179
$ws.Cells.Item(7, 4).Value = 0.1056311847747633
$ws.Cells.Item(7, 5).Value = 0.1288894075201767
$ws.Cells.Item(7, 6).Value = 2.674295448035082
$ws.Cells.Item(7, 7).Value = 2.134207411306051
$ws.Cells.Item(7, 8).Value = 1.600768445849013
$ws.Cells.Item(7, 9).Value = 3.184684259837638
$ws.Cells.Item(7, 10).Value = 0.2049917480354679
$ws.Cells.Item(8, 2).Value = 1.16825010538281
$ws.Cells.Item(8, 4).Value = 0.1057960665953246
$ws.Cells.Item(8, 5).Value = 0.129954637554043
$ws.Cells.Item(8, 6).Value = 2.834052162487836
$ws.Cells.Item(8, 7).Value = 2.318011996819791
$ws.Cells.Item(8, 8).Value = 1.675864740954125
$ws.Cells.Item(8, 9).Value = 3.795928923952317
$ws.Cells.Item(8, 10).Value = 0.212917504022812
$ws.Cells.Item(9, 2).Value = 1.539389405566055
$ws.Cells.Item(9, 4).Value = 0.1061427228871796
$ws.Cells.Item(9, 5).Value = 0.1326164825695173
$ws.Cells.Item(9, 6).Value = 3.164481609023056
$ws.Cells.Item(9, 7).Value = 2.692914156986546
$ws.Cells.Item(9, 8).Value = 1.832970800449971
$ws.Cells.Item(9, 9).Value = 4.989014804984947
$ws.Cells.Item(9, 10).Value = 0.2296922894002904
$ws.Cells.Item(10, 2).Value = 1.809644099351601
$ws.Cells.Item(10, 4).Value = 0.106412707435446
$ws.Cells.Item(10, 5).Value = 0.1349279334950673
$ws.Cells.Item(10, 6).Value = 3.418265018242977
$ws.Cells.Item(10, 7).Value = 2.977790851650184
$ws.Cells.Item(10, 8).Value = 1.954674957977772
$ws.Cells.Item(10, 9).Value = 5.861568061238813
$ws.Cells.Item(10, 10).Value = 0.2427996710749767
$ws.Cells.Item(11, 2).Value = 1.932050532542007
$ws.Cells.Item(11, 4).Value = 0.106539096441006
$ws.Cells.Item(11, 5).Value = 0.1360582219594733
$ws.Cells.Item(11, 6).Value = 3.536236937461979
$ws.Cells.Item(11, 7).Value = 3.10957700108986
$ws.Cells.Item(11, 8).Value = 2.011467405619328
$ws.Cells.Item(11, 9).Value = 6.25753114434093
$ws.Cells.Item(11, 10).Value = 0.2489397106595561
$ws.Cells.Item(12, 2).Value = 1.978324344474345
$ws.Cells.Item(12, 4).Value = 0.1065874905544781
$ws.Cells.Item(12, 5).Value = 0.1364976832432454
$ws.Cells.Item(12, 6).Value = 3.581282987014646
$ws.Cells.Item(12, 7).Value = 3.15980762891769
$ws.Cells.Item(12, 8).Value = 2.033183669984112
$ws.Cells.Item(12, 9).Value = 6.407322560669741
$ws.Cells.Item(12, 10).Value = 0.2512908683668513
$ws.Cells.Item(13, 2).Value = 1.96836198986199
$ws.Cells.Item(13, 4).Value = 0.1065770439752107
$ws.Cells.Item(13, 5).Value = 0.1364025264829394
$ws.Cells.Item(13, 6).Value = 3.571564784942183
$ws.Cells.Item(13, 7).Value = 3.148974891615353
$ws.Cells.Item(13, 8).Value = 2.02849725803037
$ws.Cells.Item(13, 9).Value = 6.375069214559915
$ws.Cells.Item(13, 10).Value = 0.2507833367248367
$ws.Cells.Item(14, 2).Value = 1.935859098657602
$ws.Cells.Item(14, 4).Value = 0.1065430670494472
$ws.Cells.Item(14, 5).Value = 0.1360941466472845
$ws.Cells.Item(14, 6).Value = 3.539935378736061
$ws.Cells.Item(14, 7).Value = 3.113702908760729
$ws.Cells.Item(14, 8).Value = 2.013249775342672
$ws.Cells.Item(14, 9).Value = 6.269857673436036
$ws.Cells.Item(14, 10).Value = 0.2491326164057455
$ws.Cells.Item(15, 2).Value = 1.915939811075134
$ws.Cells.Item(15, 4).Value = 0.1065223252554794
$ws.Cells.Item(15, 5).Value = 0.1359067491117649
$ws.Cells.Item(15, 6).Value = 3.520610273165971
$ws.Cells.Item(15, 7).Value = 3.092140582485115
$ws.Cells.Item(15, 8).Value = 2.003937769423828
$ws.Cells.Item(15, 9).Value = 6.205392518134545
$ws.Cells.Item(15, 10).Value = 0.2481249124802929
$ws.Cells.Item(16, 2).Value = 1.801633628445018
$ws.Cells.Item(16, 4).Value = 0.1064045211602469
$ws.Cells.Item(16, 5).Value = 0.1348556618023871
$ws.Cells.Item(16, 6).Value = 3.410606892642278
$ws.Cells.Item(16, 7).Value = 2.969223366258007
$ws.Cells.Item(16, 8).Value = 1.950992589575719
$ws.Cells.Item(16, 9).Value = 5.835670441680236
$ws.Cells.Item(16, 10).Value = 0.2424020215006948
$ws.Cells.Item(17, 2).Value = 1.731372246120827
$ws.Cells.Item(17, 4).Value = 0.1063331807884911
$ws.Cells.Item(17, 5).Value = 0.1342311160739662
$ws.Cells.Item(17, 6).Value = 3.343776744307405
$ws.Cells.Item(17, 7).Value = 2.894387060466613
$ws.Cells.Item(17, 8).Value = 1.918881687400301
$ws.Cells.Item(17, 9).Value = 5.608601612490133
$ws.Cells.Item(17, 10).Value = 0.2389370388622893
$ws.Cells.Item(18, 2).Value = 1.690909632146997
$ws.Cells.Item(18, 4).Value = 0.1062924826548013
$ws.Cells.Item(18, 5).Value = 0.1338793036318329
$ws.Cells.Item(18, 6).Value = 3.305574844941731
$ws.Cells.Item(18, 7).Value = 2.851549089412629
$ws.Cells.Item(18, 8).Value = 1.900546521073124
$ws.Cells.Item(18, 9).Value = 5.477907426104935
$ws.Cells.Item(18, 10).Value = 0.2369607392022033
$ws.Cells.Item(19, 2).Value = 1.67720114299982
$ws.Cells.Item(19, 4).Value = 0.1062787599065302
$ws.Cells.Item(19, 5).Value = 0.1337614551878872
$ws.Cells.Item(19, 6).Value = 3.292680760349469
$ws.Cells.Item(19, 7).Value = 2.837079934219787
$ws.Cells.Item(19, 8).Value = 1.89436145048802
$ws.Cells.Item(19, 9).Value = 5.433641513037401
$ws.Cells.Item(19, 10).Value = 0.2362944442791957
$ws.Cells.Item(20, 2).Value = 1.738856893520733
$ws.Cells.Item(20, 4).Value = 0.1063407402842653
$ws.Cells.Item(20, 5).Value = 0.1342968323355507
$ws.Cells.Item(20, 6).Value = 3.350866326720222
$ws.Cells.Item(20, 7).Value = 2.90233212842918
$ws.Cells.Item(20, 8).Value = 1.922286027038183
$ws.Cells.Item(20, 9).Value = 5.632782903627401
$ws.Cells.Item(20, 10).Value = 0.2393041633429505
$ws.Cells.Item(21, 2).Value = 1.945408141632868
$ws.Cells.Item(21, 4).Value = 0.1065530322609014
$ws.Cells.Item(21, 5).Value = 0.136184413718567
$ws.Cells.Item(21, 6).Value = 3.549215510681506
$ws.Cells.Item(21, 7).Value = 3.124054210283646
$ws.Cells.Item(21, 8).Value = 2.017722587296021
$ws.Cells.Item(21, 9).Value = 6.300765052729957
$ws.Cells.Item(21, 10).Value = 0.249616761358638
$ws.Cells.Item(22, 2).Value = 2.079940855282985
$ws.Cells.Item(22, 4).Value = 0.1066948957068323
$ws.Cells.Item(22, 5).Value = 0.1374848268051494
$ws.Cells.Item(22, 6).Value = 3.681024872415463
$ws.Cells.Item(22, 7).Value = 3.270868098814447
$ws.Cells.Item(22, 8).Value = 2.081323391753187
$ws.Cells.Item(22, 9).Value = 6.736443306619549
$ws.Cells.Item(22, 10).Value = 0.2565087896800833
$ws.Cells.Item(23, 2).Value = 2.008181005423296
$ws.Cells.Item(23, 4).Value = 0.1066188884862704
$ws.Cells.Item(23, 5).Value = 0.1367846235793699
$ws.Cells.Item(23, 6).Value = 3.610473322205905
$ws.Cells.Item(23, 7).Value = 3.192332840114375
$ws.Cells.Item(23, 8).Value = 2.047264520717476
$ws.Cells.Item(23, 9).Value = 6.503998686294722
$ws.Cells.Item(23, 10).Value = 0.2528162789516131
$ws.Cells.Item(24, 2).Value = 1.735473297351234
$ws.Cells.Item(24, 4).Value = 0.10633732165204
$ws.Cells.Item(24, 5).Value = 0.1342670994488167
$ws.Cells.Item(24, 6).Value = 3.347660443095805
$ws.Cells.Item(24, 7).Value = 2.898739583346128
$ws.Cells.Item(24, 8).Value = 1.920746533548197
$ws.Cells.Item(24, 9).Value = 5.62185100597145
$ws.Cells.Item(24, 10).Value = 0.2391381373017794
$ws.Cells.Item(25, 2).Value = 1.43940603174093
$ws.Cells.Item(25, 4).Value = 0.1060463370010964
$ws.Cells.Item(25, 5).Value = 0.1318343853988502
$ws.Cells.Item(25, 6).Value = 3.073198553213103
$ws.Cells.Item(25, 7).Value = 2.589881758271019
$ws.Cells.Item(25, 8).Value = 1.789387422476921
$ws.Cells.Item(25, 9).Value = 4.666914883151094
$ws.Cells.Item(25, 10).Value = 0.2250191059075775
